# "table module options check"
# Adds a new "action" options column (header in F1) to both worksheets,
# plus two new flag rows driven by that column:
#   sheet "one": row 7 -> on=1, action="hasMeta"; row 8 -> on=0, action="xxx"
#   sheet "two": row 7 -> on=1, action="xxx"

$wb = $excel.ActiveWorkbook

# ---- Sheet "one" ----
$ws1 = $wb.Worksheets.Item("one")
$ws1.Activate()

$ws1.Range("F1").Value = "action"

$ws1.Range("A7").Value = 1
$ws1.Range("F7").Value = "hasMeta"

$ws1.Range("A8").Value = 0
$ws1.Range("F8").Value = "xxx"

[void]$ws1.Range("A8").Select()

# ---- Sheet "two" ----
$ws2 = $wb.Worksheets.Item("two")
$ws2.Activate()

$ws2.Range("F1").Value = "action"
# match the bold/shaded header formatting used by the rest of row 1
[void]$ws2.Range("E1").Copy()
[void]$ws2.Range("F1").PasteSpecial(-4122)

$ws2.Range("A7").Value = 1
$ws2.Range("F7").Value = "xxx"

[void]$ws2.Range("A8").Select()
